$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.262.13"
$ws.Range("E2").Value = "  +11.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.60"
$ws.Range("E3").Value = "  +8.31%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.46"
$ws.Range("E5").Value = "  +4.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.573"
$ws.Range("E6").Value = "  +8.38%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.37"
$ws.Range("E8").Value = "  +7.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.70"
$ws.Range("E9").Value = "  +5.45%  "
$ws.Range("E10").Value = "  +8.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0674"
$ws.Range("E11").Value = "  +4.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0931"
$ws.Range("E12").Value = "  +3.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.077.03"
$ws.Range("E13").Value = "  +8.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.824.59"
$ws.Range("E14").Value = "  +8.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.645"
$ws.Range("E15").Value = "  +6.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.189.34"
$ws.Range("E16").Value = "  +11.54%  "
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.27"
$ws.Range("E18").Value = "  +6.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.19"
$ws.Range("E19").Value = "  +6.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "257.82"
$ws.Range("E20").Value = "  +6.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0752"
$ws.Range("E21").Value = "  +4.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.58"
$ws.Range("E23").Value = "  +6.11%  "
$ws.Range("E25").Value = "  +2.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.54"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  +5.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.117"
$ws.Range("E28").Value = "  +4.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.06"
$ws.Range("E29").Value = "  +5.53%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("E31").Value = "  +11.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0522"
$ws.Range("E32").Value = "  +5.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.21"
$ws.Range("E33").Value = "  +5.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.56"
$ws.Range("E34").Value = "  +7.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.533.51"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.79"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("E37").Value = "  +5.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.631"
$ws.Range("E38").Value = "  +5.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0189"
$ws.Range("E39").Value = "  +5.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.62"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.78"
$ws.Range("E41").Value = "  +4.30%  "
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.907"
$ws.Range("E43").Value = "  +8.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.11"
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0525"
$ws.Range("E45").Value = "  +5.31%  "
$ws.Range("E46").Value = "  +5.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.969.47"
$ws.Range("E47").Value = "  +8.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.88"
$ws.Range("E48").Value = "  +5.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.06"
$ws.Range("E49").Value = "  +16.21%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.64"
$ws.Range("E51").Value = "  +3.69%  "
